# Update the 2023 column (K) totals for the closing overview.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 4
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 7
$ws.Range("K5").Value = 18
